$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns D (open_price), E (close_price), F (high_price), G (low_price),
# H (shares_outstanding) and I (fixed_ticker) for data rows 2-43.
$ws.Cells.Item(2, 4).Value = 47.6076611892813
$ws.Cells.Item(2, 5).Value = 48.04067993164063
$ws.Cells.Item(2, 6).Value = 49.67471668521125
$ws.Cells.Item(2, 7).Value = 46.65174806777005
$ws.Cells.Item(2, 8).Value = 131684530
$ws.Cells.Item(2, 9).Value = "KLAC"
$ws.Cells.Item(3, 4).Value = 46.73817855850264
$ws.Cells.Item(3, 5).Value = 43.71403884887695
$ws.Cells.Item(3, 6).Value = 47.25730626456601
$ws.Cells.Item(3, 7).Value = 41.54687635090163
$ws.Cells.Item(3, 8).Value = 131684530
$ws.Cells.Item(3, 9).Value = "KLAC"
$ws.Cells.Item(4, 4).Value = 41.5902275540501
$ws.Cells.Item(4, 5).Value = 55.85306549072266
$ws.Cells.Item(4, 6).Value = 55.97788102380869
$ws.Cells.Item(4, 7).Value = 40.55005578696896
$ws.Cells.Item(4, 8).Value = 131684530
$ws.Cells.Item(4, 9).Value = "KLAC"
$ws.Cells.Item(5, 4).Value = 57.05376734946462
$ws.Cells.Item(5, 5).Value = 56.18156433105469
$ws.Cells.Item(5, 6).Value = 58.01822304087045
$ws.Cells.Item(5, 7).Value = 53.37206797660625
$ws.Cells.Item(5, 8).Value = 131684530
$ws.Cells.Item(5, 9).Value = "KLAC"
$ws.Cells.Item(6, 4).Value = 61.50252271170246
$ws.Cells.Item(6, 5).Value = 59.13509368896485
$ws.Cells.Item(6, 6).Value = 62.22120640698455
$ws.Cells.Item(6, 7).Value = 58.79688769979013
$ws.Cells.Item(6, 8).Value = 131684530
$ws.Cells.Item(6, 9).Value = "KLAC"
$ws.Cells.Item(7, 4).Value = 62.18629871170955
$ws.Cells.Item(7, 5).Value = 64.49485778808594
$ws.Cells.Item(7, 6).Value = 66.31785287421681
$ws.Cells.Item(7, 7).Value = 60.91701854409063
$ws.Cells.Item(7, 8).Value = 131684530
$ws.Cells.Item(7, 9).Value = "KLAC"
$ws.Cells.Item(8, 4).Value = 60.03751190356994
$ws.Cells.Item(8, 5).Value = 64.46630096435547
$ws.Cells.Item(8, 6).Value = 64.80103425985122
$ws.Cells.Item(8, 7).Value = 59.86585615396851
$ws.Cells.Item(8, 8).Value = 131684530
$ws.Cells.Item(8, 9).Value = "KLAC"
$ws.Cells.Item(9, 4).Value = 68.26765047208704
$ws.Cells.Item(9, 5).Value = 73.57553100585938
$ws.Cells.Item(9, 6).Value = 76.03928555403928
$ws.Cells.Item(9, 7).Value = 67.30808186983384
$ws.Cells.Item(9, 8).Value = 131684530
$ws.Cells.Item(9, 9).Value = "KLAC"
$ws.Cells.Item(10, 4).Value = 82.93604322982907
$ws.Cells.Item(10, 5).Value = 85.44135284423828
$ws.Cells.Item(10, 6).Value = 90.6607390252005
$ws.Cells.Item(10, 7).Value = 81.43112393092223
$ws.Cells.Item(10, 8).Value = 131684530
$ws.Cells.Item(10, 9).Value = "KLAC"
$ws.Cells.Item(11, 4).Value = 80.724171352855
$ws.Cells.Item(11, 5).Value = 81.01277923583984
$ws.Cells.Item(11, 6).Value = 88.98899418830051
$ws.Cells.Item(11, 7).Value = 79.31608943680845
$ws.Cells.Item(11, 8).Value = 131684530
$ws.Cells.Item(11, 9).Value = "KLAC"
$ws.Cells.Item(12, 4).Value = 92.92568869807278
$ws.Cells.Item(12, 5).Value = 95.8752899169922
$ws.Cells.Item(12, 6).Value = 96.65891513498076
$ws.Cells.Item(12, 7).Value = 88.96353600228096
$ws.Cells.Item(12, 8).Value = 131684530
$ws.Cells.Item(12, 9).Value = "KLAC"
$ws.Cells.Item(13, 4).Value = 93.66154412736256
$ws.Cells.Item(13, 5).Value = 97.23938751220705
$ws.Cells.Item(13, 6).Value = 107.7338006416944
$ws.Cells.Item(13, 7).Value = 92.731656597298
$ws.Cells.Item(13, 8).Value = 131684530
$ws.Cells.Item(13, 9).Value = "KLAC"
$ws.Cells.Item(14, 4).Value = 96.46596817036529
$ws.Cells.Item(14, 5).Value = 90.61441802978516
$ws.Cells.Item(14, 6).Value = 99.39619742668154
$ws.Cells.Item(14, 7).Value = 87.22996372525134
$ws.Cells.Item(14, 8).Value = 131684530
$ws.Cells.Item(14, 9).Value = "KLAC"
$ws.Cells.Item(15, 4).Value = 90.81799183270624
$ws.Cells.Item(15, 5).Value = 105.2907638549805
$ws.Cells.Item(15, 6).Value = 109.7563435315997
$ws.Cells.Item(15, 7).Value = 89.99302750410813
$ws.Cells.Item(15, 8).Value = 131684530
$ws.Cells.Item(15, 9).Value = "KLAC"
$ws.Cells.Item(16, 4).Value = 92.47323426782334
$ws.Cells.Item(16, 5).Value = 82.61760711669922
$ws.Cells.Item(16, 6).Value = 92.93352655205273
$ws.Cells.Item(16, 7).Value = 76.34501871914924
$ws.Cells.Item(16, 8).Value = 131684530
$ws.Cells.Item(16, 9).Value = "KLAC"
$ws.Cells.Item(17, 4).Value = 79.64548032681731
$ws.Cells.Item(17, 5).Value = 96.97039794921876
$ws.Cells.Item(17, 6).Value = 98.86303758391956
$ws.Cells.Item(17, 7).Value = 77.97122699453523
$ws.Cells.Item(17, 8).Value = 131684530
$ws.Cells.Item(17, 9).Value = "KLAC"
$ws.Cells.Item(18, 4).Value = 109.8451033268113
$ws.Cells.Item(18, 5).Value = 116.8089294433594
$ws.Cells.Item(18, 6).Value = 117.4778258064597
$ws.Cells.Item(18, 7).Value = 109.4144518771695
$ws.Cells.Item(18, 8).Value = 131684530
$ws.Cells.Item(18, 9).Value = "KLAC"
$ws.Cells.Item(19, 4).Value = 112.4157843135597
$ws.Cells.Item(19, 5).Value = 125.7757797241211
$ws.Cells.Item(19, 6).Value = 131.625383157238
$ws.Cells.Item(19, 7).Value = 106.9352426549819
$ws.Cells.Item(19, 8).Value = 131684530
$ws.Cells.Item(19, 9).Value = "KLAC"
$ws.Cells.Item(20, 4).Value = 149.1804100558455
$ws.Cells.Item(20, 5).Value = 156.8347320556641
$ws.Cells.Item(20, 6).Value = 162.661307766075
$ws.Cells.Item(20, 7).Value = 142.5373971159691
$ws.Cells.Item(20, 8).Value = 131684530
$ws.Cells.Item(20, 9).Value = "KLAC"
$ws.Cells.Item(21, 4).Value = 168.5383324354263
$ws.Cells.Item(21, 5).Value = 154.5167846679688
$ws.Cells.Item(21, 6).Value = 172.0064307130487
$ws.Cells.Item(21, 7).Value = 154.0599602998986
$ws.Cells.Item(21, 8).Value = 131684530
$ws.Cells.Item(21, 9).Value = "KLAC"
$ws.Cells.Item(22, 4).Value = 129.9463105284273
$ws.Cells.Item(22, 5).Value = 153.73388671875
$ws.Cells.Item(22, 6).Value = 165.9883785503893
$ws.Cells.Item(22, 7).Value = 117.6356077266978
$ws.Cells.Item(22, 8).Value = 131684530
$ws.Cells.Item(22, 9).Value = "KLAC"
$ws.Cells.Item(23, 4).Value = 183.0863794829859
$ws.Cells.Item(23, 5).Value = 188.1229553222656
$ws.Cells.Item(23, 6).Value = 197.0381562441488
$ws.Cells.Item(23, 7).Value = 176.1763970454838
$ws.Cells.Item(23, 8).Value = 131684530
$ws.Cells.Item(23, 9).Value = "KLAC"
$ws.Cells.Item(24, 4).Value = 186.8433763808394
$ws.Cells.Item(24, 5).Value = 186.4273681640625
$ws.Cells.Item(24, 6).Value = 211.3499045325026
$ws.Cells.Item(24, 7).Value = 179.8374696587072
$ws.Cells.Item(24, 8).Value = 131684530
$ws.Cells.Item(24, 9).Value = "KLAC"
$ws.Cells.Item(25, 4).Value = 248.7288907077134
$ws.Cells.Item(25, 5).Value = 265.8227233886719
$ws.Cells.Item(25, 6).Value = 301.4435546929582
$ws.Cells.Item(25, 7).Value = 244.5906997205386
$ws.Cells.Item(25, 8).Value = 131684530
$ws.Cells.Item(25, 9).Value = "KLAC"
$ws.Cells.Item(26, 4).Value = 320.730651676812
$ws.Cells.Item(26, 5).Value = 300.1258544921875
$ws.Cells.Item(26, 6).Value = 342.3252489159237
$ws.Cells.Item(26, 7).Value = 297.8131745679347
$ws.Cells.Item(26, 8).Value = 131684530
$ws.Cells.Item(26, 9).Value = "KLAC"
$ws.Cells.Item(27, 4).Value = 309.4395002674777
$ws.Cells.Item(27, 5).Value = 332.3598937988281
$ws.Cells.Item(27, 6).Value = 340.5218681926767
$ws.Cells.Item(27, 7).Value = 274.3954724267418
$ws.Cells.Item(27, 8).Value = 131684530
$ws.Cells.Item(27, 9).Value = "KLAC"
$ws.Cells.Item(28, 4).Value = 322.469731570365
$ws.Cells.Item(28, 5).Value = 356.9844970703125
$ws.Cells.Item(28, 6).Value = 372.5276113081854
$ws.Cells.Item(28, 7).Value = 304.5323901793723
$ws.Cells.Item(28, 8).Value = 131684530
$ws.Cells.Item(28, 9).Value = "KLAC"
$ws.Cells.Item(29, 4).Value = 413.9468028038758
$ws.Cells.Item(29, 5).Value = 373.7464904785156
$ws.Cells.Item(29, 6).Value = 438.8907409366319
$ws.Cells.Item(29, 7).Value = 330.3106026459549
$ws.Cells.Item(29, 8).Value = 131684530
$ws.Cells.Item(29, 9).Value = "KLAC"
$ws.Cells.Item(30, 4).Value = 355.4345476774882
$ws.Cells.Item(30, 5).Value = 307.3565368652344
$ws.Cells.Item(30, 6).Value = 356.2528322872818
$ws.Cells.Item(30, 7).Value = 300.3768500140412
$ws.Cells.Item(30, 8).Value = 131684530
$ws.Cells.Item(30, 9).Value = "KLAC"
$ws.Cells.Item(31, 4).Value = 299.0902904681363
$ws.Cells.Item(31, 5).Value = 370.4724426269531
$ws.Cells.Item(31, 6).Value = 374.9736756167053
$ws.Cells.Item(31, 7).Value = 273.1937050611585
$ws.Cells.Item(31, 8).Value = 131684530
$ws.Cells.Item(31, 9).Value = "KLAC"
$ws.Cells.Item(32, 4).Value = 296.7928520147809
$ws.Cells.Item(32, 5).Value = 306.7379760742188
$ws.Cells.Item(32, 6).Value = 324.6120440454347
$ws.Cells.Item(32, 7).Value = 242.5212125682272
$ws.Cells.Item(32, 8).Value = 131684530
$ws.Cells.Item(32, 9).Value = "KLAC"
$ws.Cells.Item(33, 4).Value = 373.49191110807
$ws.Cells.Item(33, 5).Value = 381.7398681640625
$ws.Cells.Item(33, 6).Value = 417.7078983151368
$ws.Cells.Item(33, 7).Value = 362.3844042330587
$ws.Cells.Item(33, 8).Value = 131684530
$ws.Cells.Item(33, 9).Value = "KLAC"
$ws.Cells.Item(34, 4).Value = 387.5981719239294
$ws.Cells.Item(34, 5).Value = 377.1768798828125
$ws.Cells.Item(34, 6).Value = 391.7257196820494
$ws.Cells.Item(34, 7).Value = 347.2595511126237
$ws.Cells.Item(34, 8).Value = 131684530
$ws.Cells.Item(34, 9).Value = "KLAC"
$ws.Cells.Item(35, 4).Value = 474.8880208241581
$ws.Cells.Item(35, 5).Value = 503.1929321289063
$ws.Cells.Item(35, 6).Value = 505.9539161164443
$ws.Cells.Item(35, 7).Value = 439.592536854644
$ws.Cells.Item(35, 8).Value = 131684530
$ws.Cells.Item(35, 9).Value = "KLAC"
$ws.Cells.Item(36, 4).Value = 450.7142458886467
$ws.Cells.Item(36, 5).Value = 461.1206359863281
$ws.Cells.Item(36, 6).Value = 497.6607896764617
$ws.Cells.Item(36, 7).Value = 440.1213238677309
$ws.Cells.Item(36, 8).Value = 131684530
$ws.Cells.Item(36, 9).Value = "KLAC"
$ws.Cells.Item(37, 4).Value = 564.2889661613768
$ws.Cells.Item(37, 5).Value = 584.7947387695312
$ws.Cells.Item(37, 6).Value = 648.5272193608441
$ws.Cells.Item(37, 7).Value = 533.9682686077467
$ws.Cells.Item(37, 8).Value = 131684530
$ws.Cells.Item(37, 9).Value = "KLAC"
$ws.Cells.Item(38, 4).Value = 690.63720549312
$ws.Cells.Item(38, 5).Value = 680.0704345703125
$ws.Cells.Item(38, 6).Value = 708.5740331181013
$ws.Cells.Item(38, 7).Value = 614.8348222060938
$ws.Cells.Item(38, 8).Value = 131684530
$ws.Cells.Item(38, 9).Value = "KLAC"
$ws.Cells.Item(39, 4).Value = 813.8317657473086
$ws.Cells.Item(39, 5).Value = 813.713134765625
$ws.Cells.Item(39, 6).Value = 886.1304098347954
$ws.Cells.Item(39, 7).Value = 740.8608593165171
$ws.Cells.Item(39, 8).Value = 131684530
$ws.Cells.Item(39, 9).Value = "KLAC"
$ws.Cells.Item(40, 4).Value = 771.0283893441335
$ws.Cells.Item(40, 5).Value = 659.861328125
$ws.Cells.Item(40, 6).Value = 825.0371529984118
$ws.Cells.Item(40, 7).Value = 648.3920676505629
$ws.Cells.Item(40, 8).Value = 131684530
$ws.Cells.Item(40, 9).Value = "KLAC"
$ws.Cells.Item(41, 4).Value = 630.3668796758534
$ws.Cells.Item(41, 5).Value = 733.2115478515625
$ws.Cells.Item(41, 6).Value = 781.7685534175789
$ws.Cells.Item(41, 7).Value = 627.5859636989279
$ws.Cells.Item(41, 8).Value = 131684530
$ws.Cells.Item(41, 9).Value = "KLAC"
$ws.Cells.Item(42, 4).Value = 670.9212736039041
$ws.Cells.Item(42, 5).Value = 699.480224609375
$ws.Cells.Item(42, 6).Value = 718.5427747123849
$ws.Cells.Item(42, 7).Value = 548.8116279554603
$ws.Cells.Item(42, 8).Value = 131684530
$ws.Cells.Item(42, 9).Value = "KLAC"
$ws.Cells.Item(43, 4).Value = 886.8499819857421
$ws.Cells.Item(43, 5).Value = 877.1211547851562
$ws.Cells.Item(43, 6).Value = 943.8159729963244
$ws.Cells.Item(43, 7).Value = 858.1324505132194
$ws.Cells.Item(43, 8).Value = 131684530
$ws.Cells.Item(43, 9).Value = "KLAC"
